$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text without leaving a residual
# NumberFormat style behind (Excel auto-coerces clean numeric-looking
# strings like "1.001" or "306.05" to doubles otherwise).
function Set-TextValue($cell, $text) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.ClearFormats()
}

$ws.Range("D2").Value = "23.447.83"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").Value = "1.638.37"
$ws.Range("E3").Value = "  +2.29%  "
$ws.Range("E4").Value = "  +0.04%  "
Set-TextValue "D5" "1.001"
$ws.Range("E5").Value = "  +0.00%  "
Set-TextValue "D6" "306.05"
$ws.Range("E6").Value = "  +0.88%  "
Set-TextValue "D7" "0.3763"
$ws.Range("E7").Value = "  -0.53%  "
Set-TextValue "D8" "52.09"
$ws.Range("E8").Value = "  +0.10%  "
Set-TextValue "D9" "0.3640"
$ws.Range("E9").Value = "  +0.68%  "
Set-TextValue "D10" "1.261"
$ws.Range("E10").Value = "  -0.73%  "
Set-TextValue "D11" "0.08143"
$ws.Range("E11").Value = "  +0.27%  "
Set-TextValue "D12" "1.001"
$ws.Range("E12").Value = "  +0.07%  "
Set-TextValue "D13" "22.92"
$ws.Range("E13").Value = "  +0.81%  "
Set-TextValue "D14" "6.622"
$ws.Range("E14").Value = "  +0.62%  "
Set-TextValue "D15" "0.00001274"
$ws.Range("E15").Value = "  +2.52%  "
Set-TextValue "D16" "7.355"
$ws.Range("E16").Value = "  -0.73%  "
$ws.Range("D17").Value = "1.637.85"
$ws.Range("E17").Value = "  +2.25%  "
Set-TextValue "D18" "94.63"
$ws.Range("E18").Value = "  +0.45%  "
Set-TextValue "D19" "0.06911"
$ws.Range("E19").Value = "  +0.41%  "
Set-TextValue "D20" "18.16"
$ws.Range("E20").Value = "  +0.39%  "
Set-TextValue "D21" "6.539"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D23").Value = "23.449.49"
$ws.Range("E23").Value = "  +1.15%  "
Set-TextValue "D24" "12.77"
$ws.Range("E24").Value = "  -1.63%  "
Set-TextValue "D25" "3.064"
$ws.Range("E25").Value = "  +2.70%  "
Set-TextValue "D26" "2.419"
$ws.Range("E26").Value = "  +0.88%  "
Set-TextValue "D27" "21.21"
$ws.Range("E27").Value = "  -0.15%  "
Set-TextValue "D28" "150.73"
$ws.Range("E28").Value = "  +0.87%  "
Set-TextValue "D29" "5.355"
Set-TextValue "D30" "136.96"
$ws.Range("E30").Value = "  +2.05%  "
Set-TextValue "D31" "2.306"
$ws.Range("E31").Value = "  -3.28%  "
$ws.Range("D32").Value = "1.819.70"
$ws.Range("E32").Value = "  +2.34%  "
Set-TextValue "D33" "6.769"
$ws.Range("E33").Value = "  -0.05%  "
Set-TextValue "D34" "0.9647"
$ws.Range("E34").Value = "  -0.43%  "
Set-TextValue "D35" "0.02838"
$ws.Range("E35").Value = "  +4.61%  "
Set-TextValue "D36" "10.30"
$ws.Range("E36").Value = "  +0.02%  "
Set-TextValue "D37" "0.07300"
$ws.Range("E37").Value = "  -2.76%  "
Set-TextValue "D38" "0.2528"
$ws.Range("E38").Value = "  +1.00%  "
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D39" "0.08831"
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D40" "6.113"
$ws.Range("E40").Value = "  +0.33%  "
Set-TextValue "D41" "1.377"
$ws.Range("E41").Value = "  +1.28%  "
Set-TextValue "D42" "0.7089"
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D43" "16.32"
$ws.Range("E43").Value = "  +4.31%  "
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D44" "12.48"
$ws.Range("E44").Value = "  -0.22%  "
Set-TextValue "D45" "0.6545"
$ws.Range("E45").Value = "  +0.27%  "
Set-TextValue "D46" "2.336"
$ws.Range("E46").Value = "  +0.96%  "
$ws.Range("E47").Value = "  +0.03%  "
Set-TextValue "D48" "4.015"
$ws.Range("E48").Value = "  -0.11%  "
Set-TextValue "D49" "0.07970"
$ws.Range("E49").Value = "  -0.01%  "
Set-TextValue "D50" "128.91"
$ws.Range("E50").Value = "  -2.41%  "
Set-TextValue "D51" "1.202"
$ws.Range("E51").Value = "  +0.26%  "
